$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of (row, date serial, B value, C value) appended after the last existing row (97)
$rows = @(
    @(98,  45658, 0.231053709604005, 0.0884910594688836),
    @(99,  45689, 0.24505006827492,  0.0925520379308145),
    @(100, 45717, 0.203326505841299, 0.124402015606064)
)

foreach ($r in $rows) {
    $rowIndex = $r[0]
    $dateVal  = $r[1]
    $bVal     = $r[2]
    $cVal     = $r[3]

    # Carry the date-formatted style from column A of the row above onto the
    # new row's A cell (keeps the same style index instead of minting a new one)
    $ws.Range("A" + ($rowIndex - 1)).Copy()
    $ws.Range("A" + $rowIndex).PasteSpecial(-4122)

    $ws.Cells.Item($rowIndex, 1).Value = $dateVal
    $ws.Cells.Item($rowIndex, 2).Value = $bVal
    $ws.Cells.Item($rowIndex, 3).Value = $cVal
}
